# ---------------------------------------------------------------------------
# Adds a new "ODI Bowling Extra" worksheet (scraped bowling-extra attributes)
# and cleans up stray empty cells left behind in "ODI Batting Extra".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "ODI Batting Extra" (sheet 4): drop the empty placeholder cells that
#    carry no value (they were written out as blank inlineStr cells).
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCellsToClear = @(
    "E3",
    "B4", "C4", "D4", "E4",
    "C6", "D6", "E6",
    "B7", "C7", "D7", "E7",
    "B8", "C8", "D8", "E8",
    "B9", "C9", "D9", "E9",
    "C10", "D10", "E10",
    "B11", "C11", "D11", "E11",
    "B13", "C13", "D13", "E13",
    "B14", "C14", "D14", "E14",
    "E16",
    "B21", "C21", "D21", "E21", "F21"
)

foreach ($ref in $emptyCellsToClear) {
    $battingExtra.Range($ref).ClearContents()
}

# ---------------------------------------------------------------------------
# 2) Add the new "ODI Bowling Extra" worksheet straight after
#    "ODI Batting Extra" and populate it with the scraped data.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row (bold, centred, thin border - same look as the other sheets).
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $bowlingExtra.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @{A="4431"; B="0"; C="10.00%"},
    @{A="4454"; B="0"; C=$null},
    @{A="4456"; B=$null; C=$null},
    @{A="4457"; B="0"; C="20.00%"},
    @{A="4469"; B="0"; C=$null},
    @{A="4470"; B="0"; C=$null},
    @{A="4471"; B=$null; C=$null},
    @{A="4598"; B=$null; C=$null},
    @{A="4599"; B=$null; C=$null},
    @{A="4602"; B="0"; C="10.00%"},
    @{A="4619"; B=$null; C=$null},
    @{A="4620"; B="1"; C="30.00%"},
    @{A="4622"; B=$null; C=$null},
    @{A="4663"; B=$null; C=$null},
    @{A="4698"; B="0"; C="10.00%"},
    @{A="4699"; B="0"; C="20.00%"},
    @{A="4700"; B="0"; C="30.00%"},
    @{A="4711"; B="0"; C="20.00%"},
    @{A="4713"; B="0"; C="40.00%"},
    @{A="4717"; B="0"; C="20.00%"}
)

$r = 2
foreach ($row in $rows) {
    if ($null -ne $row.A) {
        $cell = $bowlingExtra.Cells.Item($r, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $row.A
    }
    if ($null -ne $row.B) {
        $cell = $bowlingExtra.Cells.Item($r, 2)
        $cell.NumberFormat = "@"
        $cell.Value = $row.B
    }
    if ($null -ne $row.C) {
        $cell = $bowlingExtra.Cells.Item($r, 3)
        $cell.NumberFormat = "@"
        $cell.Value = $row.C
    }
    $r++
}

Write-Host "Added sheet 'ODI Bowling Extra' and cleaned empty cells in 'ODI Batting Extra'."
